# DE_table2_F6_dim10.xlsx edit
#
# Summary of the change:
#  - The "Run 50" column (last run, previously column AZ) is removed
#    entirely -- the sheet now only reports Run 0 .. Run 49 (50 runs).
#  - The trailing "Mean" column (previously BA, now shifted left into AZ
#    once "Run 50" is deleted) is recomputed over the remaining 50 runs.
#  - The first column header "Gen" becomes "MaxFES", and its values switch
#    from generation counts to the normalised 0..1 MaxFES fractions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Run 50" column (AZ). Excel shifts everything after it
# (just the "Mean" column, BA) one column to the left, so "Mean" becomes
# the new column AZ -- matching dimension A1:AZ14 / spans 1:52.
$ws.Columns("AZ").Delete()

# Rename header "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# New MaxFES values (column A, rows 2-14)
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 0.001
$ws.Range("A4").Value = 0.01
$ws.Range("A5").Value = 0.1
$ws.Range("A6").Value = 0.2
$ws.Range("A7").Value = 0.3
$ws.Range("A8").Value = 0.4
$ws.Range("A9").Value = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# Recomputed Mean values (now column AZ), averaged over Run 0 .. Run 49 only
$ws.Range("AZ2").Value = 13.30211237
$ws.Range("AZ3").Value = 10.93647393
$ws.Range("AZ4").Value = 5.45155992
$ws.Range("AZ5").Value = 5.30507009
$ws.Range("AZ6").Value = 5.30507009
$ws.Range("AZ7").Value = 5.30507009
$ws.Range("AZ8").Value = 5.30507009
$ws.Range("AZ9").Value = 5.30507009
$ws.Range("AZ10").Value = 5.30507009
$ws.Range("AZ11").Value = 5.30507009
$ws.Range("AZ12").Value = 5.30507009
$ws.Range("AZ13").Value = 5.30507009
$ws.Range("AZ14").Value = 5.30507009
